$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1 & 2) "GIT version control and GitHub" -> "Git version control and GitHub"
#    Split into three runs: "G" | "it" | " version control and GitHub"
#    with the _GoBack bookmark sitting between run 2 and run 3.
#    Adding a new "_GoBack" bookmark automatically relocates Word's single
#    hidden _GoBack bookmark away from the phone-number paragraph (its
#    original location), which is exactly the other half of this change.
# ---------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("GIT version control and GitHub", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $gitStart = $rng.Start

    # Fix the casing in place first (same formatting, so this alone would
    # still coalesce into a single run on save).
    $rG = $d.Range($gitStart, $gitStart + 1)
    $rG.Text = "G"
    $rIt = $d.Range($gitStart + 1, $gitStart + 3)
    $rIt.Text = "it"

    # Drop a transient bookmark between "G" and "it" just to force the
    # run boundary; it gets removed again right after.
    $tempSplit = $d.Range($gitStart + 1, $gitStart + 1)
    $d.Bookmarks.Add("ZZ_TEMP_SPLIT", $tempSplit)

    # Re-home the real _GoBack bookmark between "it" and " version control...".
    $goBackRange = $d.Range($gitStart + 3, $gitStart + 3)
    $d.Bookmarks.Add("_GoBack", $goBackRange)

    $d.Bookmarks("ZZ_TEMP_SPLIT").Delete()
}

# ---------------------------------------------------------------------
# 3) Insert a new bullet "Documentation languages, such as Markdown"
#    right before the "Parallel / concurrent programming" bullet.
# ---------------------------------------------------------------------
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("Parallel / concurrent programming", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $targetPara = $rng2.Paragraphs(1)
    $insertPos = $targetPara.Range.Start
    $targetPara.Range.InsertParagraphBefore()
    $newRange = $d.Range($insertPos, $insertPos)
    $newRange.Text = "Documentation languages, such as Markdown"
}

# ---------------------------------------------------------------------
# 4) "Unit Testing" -> "Unit testing"
#    Split into three runs: "Unit " | "t" | "esting"
# ---------------------------------------------------------------------
$rng3 = $d.Content
$found3 = $rng3.Find.Execute("Unit Testing", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found3) {
    $utStart = $rng3.Start

    $rT = $d.Range($utStart + 5, $utStart + 6)
    $rT.Text = "t"

    $bmA = $d.Range($utStart + 5, $utStart + 5)
    $d.Bookmarks.Add("ZZ_TEMP_A", $bmA)
    $bmB = $d.Range($utStart + 6, $utStart + 6)
    $d.Bookmarks.Add("ZZ_TEMP_B", $bmB)

    $d.Bookmarks("ZZ_TEMP_A").Delete()
    $d.Bookmarks("ZZ_TEMP_B").Delete()
}
